# Add team record (Wins/Losses/Ties) columns to the roster sheet.
# New columns: AD = Wins, AE = Losses, AF = Ties
# Every player row (2-51) gets the same team record: 56 wins, 106 losses, 0 ties.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - set text values for the three new columns.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header style (bold font, thin border, centered alignment) from an
# existing header cell (A1) onto the new header cells so they match the rest
# of row 1.
$headerStyleSource = $ws.Range("A1")
$newHeaderCells = $ws.Range("AD1:AF1")
$headerStyleSource.Copy()
$newHeaderCells.PasteSpecial(-4122)

# Data rows (2-51) - fill in the team record values.
for ($r = 2; $r -le 51; $r++) {
    $ws.Cells.Item($r, 30).Value = 56
    $ws.Cells.Item($r, 31).Value = 106
    $ws.Cells.Item($r, 32).Value = 0
}
